# Insert a new weekly price record as row 359 (pushing existing rows 359-375
# down to 360-376), matching the new "Haba" observation added to the
# Vega Central Mapocho de Santiago consolidated sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 359:375 down to 360:376 by inserting a fresh row at 359.
# Excel's EntireRow.Insert copies the formatting (incl. the date number
# format on column D) down from the row above automatically.
$ws.Rows(359).Insert()

# Populate the newly inserted row with the new record's values.
$ws.Cells.Item(359, 1).Value = 9
$ws.Cells.Item(359, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(359, 3).Value = "Metropolitana"
$ws.Cells.Item(359, 4).Value = 45147
$ws.Cells.Item(359, 5).Value = 13
$ws.Cells.Item(359, 6).Value = 100112026
$ws.Cells.Item(359, 7).Value = "Haba"
$ws.Cells.Item(359, 8).Value = "Sin especificar"
$ws.Cells.Item(359, 9).Value = "Primera"
$ws.Cells.Item(359, 10).Value = 52
$ws.Cells.Item(359, 11).Value = 10000
$ws.Cells.Item(359, 12).Value = 13000
$ws.Cells.Item(359, 13).Value = 11500
$ws.Cells.Item(359, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(359, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(359, 16).Value = 460
$ws.Cells.Item(359, 17).Value = 25
$ws.Cells.Item(359, 18).Value = "Hortaliza"
